# Añadiendo lectura de S3
# Adds a new booking-source row (BODEGAS FERRERA S.L / Viator) and a new
# "Emails_Reservas" entry for an existing row (Viator booking channel).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New reservation-inbox email for row 6 (Cocodrilo Park / VIATOR,INC channel)
$ws.Range("Q6").Value2 = "booking@t1.viator.com"

# New activity row: "Visita guiada a viñedo ecológico..." (row 34) gains its
# provider / fiscal data, pricing and reservation-email content.
$ws.Range("D34").Value2 = "BODEGAS FERRERA S.L"
$ws.Range("E34").Value2 = "BODEGAS FERRERA S.L"
$ws.Range("F34").Value2 = "B76684992"

$ws.Range("H34").Value2 = 45
$ws.Range("J34").Value2 = 45
$ws.Range("L34").Value2 = 22
$ws.Range("M34").Value2 = 0.28039999999999998

$ws.Range("N34").Value2 = "3786 "
$ws.Range("O34").Value2 = "visitas@bodegasferrera.com"
$ws.Range("P34").Value2 = "Date of the activity, Time of the activity, Total adults, Total children (If applicable, leave blank if not.) , Name of the primary client, Phone number"
